$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

# Add ORCID user IDs for row 4 (E. Taylor / Crockford) and row 5 (Emily Peacock)
$ws.Range("F4").Value = "0000-0002-2122-0462"
$ws.Range("F5").Value = "0000-0003-0194-7282"

# Set the new column F width to match the corrected layout
# (engine quantizes ColumnWidth to 1/6-character pixel steps; 25.33 is the
# input that lands closest to the authored width of 26.21875 characters)
$ws.Range("F1").ColumnWidth = 25.33

# Update the selected cell in the sheet view
$ws.Range("D20").Select()
